$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2025-06-13 Friday"; New = "2025-06-14 Saturday" },
    @{ Old = "563÷9="; New = "535÷4=" },
    @{ Old = "105÷3="; New = "697÷5=" },
    @{ Old = "710÷9="; New = "330÷4=" },
    @{ Old = "856÷2="; New = "841÷8=" },
    @{ Old = "849÷4="; New = "901÷9=" },
    @{ Old = "112÷5="; New = "356÷6=" },
    @{ Old = "992÷7="; New = "209÷3=" },
    @{ Old = "366÷9="; New = "869÷7=" },
    @{ Old = "211÷7="; New = "212÷7=" },
    @{ Old = "969÷7="; New = "791÷5=" },
    @{ Old = "968÷9="; New = "761÷9=" },
    @{ Old = "260÷3="; New = "738÷5=" },
    @{ Old = "420÷6="; New = "833÷6=" },
    @{ Old = "501÷4="; New = "956÷8=" },
    @{ Old = "139÷5="; New = "540÷9=" },
    @{ Old = "525÷7="; New = "102÷2=" },
    @{ Old = "577÷3="; New = "701÷6=" },
    @{ Old = "336÷8="; New = "501÷5=" },
    @{ Old = "161÷3="; New = "721÷4=" },
    @{ Old = "608÷4="; New = "399÷4=" },
    @{ Old = "408÷4="; New = "810÷3=" },
    @{ Old = "126÷6="; New = "441÷3=" },
    @{ Old = "522÷3="; New = "831÷3=" },
    @{ Old = "242÷9="; New = "719÷8=" },
    @{ Old = "292÷6="; New = "849÷3=" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.New, 2)
}
